# feat: add 2022-Q1 data
#
# The workbook tracks one "snapshot" worksheet per quarter (e.g. "2020-Q4",
# "2021-Q4", ...) plus a rolling "总计" (totals) sheet that lists every
# snapshot taken so far, newest first.
#
# Adding a new quarter means:
#   1. The current "总计" sheet's data actually belongs to the new quarter
#      snapshot, so it gets renamed to "2022-Q1" and its single data row is
#      replaced with that quarter's fund-holding figures.
#   2. A brand new "总计" sheet is appended at the end, reusing the old
#      totals table but with a new first row for "2022-Q1" prepended above
#      the existing history.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Turn the old "总计" sheet into the new "2022-Q1" snapshot sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("B2").Value = "004250"
$q1.Range("C2").Value = "银河量化优选混合"
$q1.Range("D2").Value = "0.39"
$q1.Range("E2").Value = "80.03"
$q1.Range("F2").Value = "1.89"
$q1.Range("G2").Value = "0.0074"
$q1.Range("H2").Value = 4

# ---------------------------------------------------------------------
# 2. Append a fresh "总计" sheet after the last existing sheet, rebuilding
#    the running history table with "2022-Q1" as the newest entry.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.01

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2020-Q4"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.24
